$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update ledStatus for the Z005 project rows (E3 and E6) from numeric 9 to text "Concluido"
$ws.Range("E3").Value = "Concluido"
$ws.Range("E6").Value = "Concluido"
